$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.67"
$ws.Range("E2").Value = "'-9.61%"
$ws.Range("D3").Value = "'40.39"
$ws.Range("E3").Value = "'-2.41%"
$ws.Range("D4").Value = "'5.031"
$ws.Range("E4").Value = "'-4.45%"
$ws.Range("D5").Value = "'0.07292"
$ws.Range("E5").Value = "'-5.87%"
$ws.Range("D6").Value = "'4.283"
$ws.Range("E6").Value = "'-1.29%"
$ws.Range("D7").Value = "'1.519"
$ws.Range("E7").Value = "'-13.56%"
$ws.Range("D8").Value = "'0.9182"
$ws.Range("E8").Value = "'-2.74%"
$ws.Range("D9").Value = "'0.1178"
$ws.Range("E9").Value = "'-5.13%"
$ws.Range("D10").Value = "'0.1729"
$ws.Range("E10").Value = "'-8.44%"
$ws.Range("D11").Value = "'0.08647"
$ws.Range("D12").Value = "'0.04176"
$ws.Range("E12").Value = "'-3.03%"
$ws.Range("D13").Value = "'0.1054"
$ws.Range("E13").Value = "'0.23%"
$ws.Range("D14").Value = "'0.001269"
$ws.Range("E14").Value = "'-0.93%"
$ws.Range("D15").Value = "'0.005835"
$ws.Range("E15").Value = "'-0.28%"
$ws.Range("D16").Value = "'3.400"
$ws.Range("E16").Value = "'1.75%"
$ws.Range("D18").Value = "'0.3289"
$ws.Range("D19").Value = "'7.868"
$ws.Range("E19").Value = "'0.94%"
$ws.Range("D20").Value = "'0.1342"
$ws.Range("E20").Value = "'-0.07%"
$ws.Range("D22").Value = "'0.03870"
$ws.Range("E22").Value = "'-4.20%"
$ws.Range("D23").Value = "'0.001269"
$ws.Range("E23").Value = "'0.11%"
$ws.Range("D24").Value = "'0.003821"
$ws.Range("E24").Value = "'-7.40%"
$ws.Range("D25").Value = "'0.0001283"
$ws.Range("E25").Value = "'0.84%"
$ws.Range("D26").Value = "'0.0003727"
$ws.Range("E26").Value = "'-95.02%"
$ws.Range("D38").Value = "'0.02319"
$ws.Range("E38").Value = "'-9.34%"
$ws.Range("D39").Value = "'0.04969"
$ws.Range("E39").Value = "'-6.91%"
$ws.Range("D40").Value = "'0.006935"
$ws.Range("E40").Value = "'248.35%"
$ws.Range("D41").Value = "'0.007687"
$ws.Range("E41").Value = "'-0.82%"
$ws.Range("D42").Value = "'0.1275"
$ws.Range("E42").Value = "'-3.24%"
$ws.Range("D43").Value = "'0.007354"
$ws.Range("E43").Value = "'4.37%"
$ws.Range("E44").Value = "'-14.47%"
$ws.Range("D45").Value = "'0.3121"
$ws.Range("E45").Value = "'-1.51%"
$ws.Range("D46").Value = "'0.00006434"
$ws.Range("E46").Value = "'-3.65%"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("D48").Value = "'0.008538"
$ws.Range("E48").Value = "'-95.75%"
$ws.Range("E49").Value = "'0.01%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'-0.01%"
